$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix product types (column F) for several rows
$ws.Range("F2").Value = "scenario"
$ws.Range("F7").Value = "supplement"
$ws.Range("F12").Value = "supplement"
$ws.Range("F16").Value = "supplement"
$ws.Range("F19").Value = "scenario"
$ws.Range("F20").Value = "supplement"

# Update the active selection to reflect the last edited cell
$ws.Range("F21").Select()
